# Implemented data driven approach
# Replace the hard-coded single-user login creds with three data-driven
# users (user1/2/3@gmail.com) and add a third row of test data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old hyperlinks before rewriting the cell values/rel ids.
$ws.Hyperlinks.Delete()

# Row 2 -> user1
$ws.Range("A2").Value = "user1@gmail.com"
$ws.Range("B2").Value = "ThisIsUser1Pass"

# Row 3 -> user2
$ws.Range("A3").Value = "user2@gmail.com"
$ws.Range("B3").Value = "ThisIsUser2Pass"
$ws.Range("B3").Style = "Normal"

# Row 4 (new) -> user3. Password entered before the email, same order the
# original author typed the data in, so shared-string indices line up.
$ws.Range("B4").Value = "ThisIsUser3Pass"
$ws.Range("A4").Value = "user3@gmail.com"

# Re-create the mailto hyperlinks, matching the order from the diff
# (A3, then A2, then A4).
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:user2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:user1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:user3@gmail.com")

# Keep the "Hyperlink" cell style on the e-mail column.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"

# Page setup tweak that shipped with the same commit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leftover selection state from the author's last click.
$ws.Range("B12").Select()
